$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.316.06'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.903.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5135'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3943'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08467'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.120'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.266'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.902.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.362'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06734'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.044'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.311.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.222'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.120.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.460'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.062'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1050'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.096'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.655'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02483'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.26%  '
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.138'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2201'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('E39').Value = '  +4.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.137'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6523'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.235'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6059'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.681'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.060'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.86'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.48%  '
